# Edit: remove the "Cohort" line from the StatQuery (B2) query text,
# update selection to B2, and adjust row heights to match the
# newer Excel build's text-wrapping metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Airedale Terrier', 'Bluetick Hound', 'Norfolk Terrier', 'Norwegian Elkhound', 'Scottish Terrier', 'Welsh Springer Spaniel', 'Wheaten Terrier']and diag.disease_term in ['Bladder Cancer'] and demo.sex in ['Male', 'Female'] and demo.neutered_indicator IN ['No', 'Yes']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newQuery

$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 290
$ws.Rows.Item(4).RowHeight = 275.5

$ws.Range("B2").Select()
